$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.669.78"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3
$ws.Range("D3").Value = "3.360.64"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "558.97"
$ws.Range("E5").Value = "  -0.54%  "

# Row 6
$ws.Range("D6").Value = "153.75"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").Value = "3.358.28"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  -2.41%  "

# Row 11
$ws.Range("D11").Value = "0.120"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
$ws.Range("D12").Value = "0.429"
$ws.Range("E12").Value = "  -1.94%  "

# Row 13
$ws.Range("D13").Value = "3.937.31"
$ws.Range("E13").Value = "  -0.12%  "

# Row 14
$ws.Range("E14").Value = "  -3.90%  "

# Row 15
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  +2.14%  "

# Row 16
$ws.Range("D16").Value = "26.69"
$ws.Range("E16").Value = "  -1.93%  "

# Row 17
$ws.Range("D17").Value = "62.721.91"
$ws.Range("E17").Value = "  -0.43%  "

# Row 18
$ws.Range("D18").Value = "3.371.07"
$ws.Range("E18").Value = "  +0.29%  "

# Row 19
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  -4.52%  "

# Row 20
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21
$ws.Range("D21").Value = "371.93"
$ws.Range("E21").Value = "  -4.50%  "

# Row 22
$ws.Range("D22").Value = "7.94"
$ws.Range("E22").Value = "  -5.95%  "

# Row 23
$ws.Range("E23").Value = "  -0.82%  "

# Row 24
$ws.Range("D24").Value = "70.71"
$ws.Range("E24").Value = "  +0.33%  "

# Row 25
$ws.Range("D25").Value = "0.523"
$ws.Range("E25").Value = "  -3.42%  "

# Row 26
$ws.Range("D26").Value = "0.0000113"
$ws.Range("E26").Value = "  +15.97%  "

# Row 27
$ws.Range("D27").Value = "9.36"
$ws.Range("E27").Value = "  +5.90%  "

# Row 28
$ws.Range("D28").Value = "0.175"
$ws.Range("E28").Value = "  -2.94%  "

# Row 29
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").Value = "6.01"
$ws.Range("E30").Value = "  +5.96%  "

# Row 31
$ws.Range("D31").Value = "6.47"
$ws.Range("E31").Value = "  -2.04%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.33"
$ws.Range("E32").Value = "  +1.59%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.96"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "22.94"
$ws.Range("E34").Value = "  -0.51%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").Value = "6.67"
$ws.Range("E36").Value = "  -0.86%  "

# Row 37
$ws.Range("D37").Value = "158.91"
$ws.Range("E37").Value = "  -1.25%  "

# Row 38
$ws.Range("E38").Value = "  -3.50%  "

# Row 39
$ws.Range("D39").Value = "0.0759"
$ws.Range("E39").Value = "  +1.42%  "

# Row 40
$ws.Range("D40").Value = "2.895.88"
$ws.Range("E40").Value = "  +1.97%  "

# Row 41
$ws.Range("D41").Value = "26.74"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("D42").Value = "1.80"
$ws.Range("E42").Value = "  -5.18%  "

# Row 43
$ws.Range("D43").Value = "0.0314"
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("D44").Value = "40.89"
$ws.Range("E44").Value = "  +0.19%  "

# Row 45
$ws.Range("D45").Value = "4.27"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
$ws.Range("E46").Value = "  -1.71%  "

# Row 47
$ws.Range("D47").Value = "22.74"
$ws.Range("E47").Value = "  +2.06%  "

# Row 48
$ws.Range("D48").Value = "1.05"
$ws.Range("E48").Value = "  -0.27%  "

# Row 49
$ws.Range("E49").Value = "  +13.85%  "

# Row 50
$ws.Range("D50").Value = "6.30"
$ws.Range("E50").Value = "  -0.35%  "

# Row 51
$ws.Range("D51").Value = "0.821"
$ws.Range("E51").Value = "  +1.37%  "
